$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Insert a new row above the existing "ACCULT" row so the new "CE Part 1"
# section header lands at row 2 and the original rows (ACCULT/MACV/MEIM)
# shift down to rows 3-5.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "CE Part 1"
$ws.Range("B2").Value = "Culture-related Measures"

$ws.Range("A6").Value = "VIA"
$ws.Range("B6").Value = "Vancouver Index of Acculturation"

$ws.Range("A7").Value = "CE Part 2"
$ws.Range("B7").Value = "Family-related Measures"

$ws.Range("A8").Value = "PBI"
$ws.Range("B8").Value = "children's report of Parental Behavioral Inventory"

$ws.Range("A9").Value = "MNBS"
$ws.Range("B9").Value = "Multidimensional Neglecful Behavior Scale"

$ws.Range("B10").Value = "Parental Monitoring"
$ws.Range("A10").Value = "PM"

$ws.Range("A11").Value = "FES"
$ws.Range("B11").Value = "Family Environment Scale"

$ws.Range("A12").Value = "PET"
$ws.Range("B12").Value = "Pet Ownership"

$ws.Range("A13").Value = "CE Part 3"
$ws.Range("B13").Value = "Peers- and Family-related Measures"

$ws.Range("A14").Value = "PBP"
$ws.Range("B14").Value = "Peer Behavior Profile"

$ws.Range("B15").Value = "Peer Network Health: protective scale"
$ws.Range("A15").Value = "PNH"

$ws.Range("A16").Value = "RPI"
$ws.Range("B16").Value = "Resistance to Peer Influence"

$ws.Range("A17").Value = "SAG"
$ws.Range("B17").Value = "School Attendance and Grades"

$ws.Range("A18").Value = "SRPF"
$ws.Range("B18").Value = "School Risk and Protective Factors scale"

$ws.Range("A19").Value = "CE Part 4"
$ws.Range("B19").Value = "Other CE-related Measures"

$ws.Range("A20").Value = "PSB"
$ws.Range("B20").Value = "ProSocial Behavior scale"

$ws.Range("A21").Value = "WPS"
$ws.Range("B21").Value = "Wills Problem Solving scale"

# Section-header rows: centered horizontal alignment (matches rows 2/7/13/19)
$headerRows = @(2,7,13,19)
foreach ($r in $headerRows) {
    $ws.Range("A" + $r + ":B" + $r).HorizontalAlignment = -4108
}

$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Range("C26").Select()
